$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new departure record as row 28 (mirrors the pattern of the
# existing rows: Number, Date, Time, Flight, To, Short, Airline, Model,
# Aircraft ID, Status, <blank>, Difference, <blank>).
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = "Monday, Jan 16"
$ws.Range("C28").Value = "8:25 AM"
$ws.Range("D28").Value = "FR6112"
$ws.Range("E28").Value = "Gdansk"
$ws.Range("F28").Value = "(GDN)"
$ws.Range("G28").Value = "Ryanair "
$ws.Range("H28").Value = "B738"
$ws.Range("I28").Value = "(SP-RSW)"
$ws.Range("J28").Value = "8:28 AM"
$ws.Range("L28").Value = "0 hours, 3 minutes"

# K28 and M28 stay blank, same as every other row, but Excel still
# materialises an (empty) cell entry for them. Touch each with a
# default-value property write so the cell gets created without
# introducing a new style.
$ws.Range("K28").Font.Bold = $false
$ws.Range("M28").Font.Bold = $false
